$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty Hours/Activities cells for rows 17-21
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = "Meetings, study, web sockets research and oanic buttons"

$ws.Range("E18").Value = 6
$ws.Range("F18").Value = "Meetings, study and panic buttons"

$ws.Range("E19").Value = 4
$ws.Range("F19").Value = "Meetings and study"

$ws.Range("E20").Value = 6
$ws.Range("F20").Value = "Study and general fixes "

$ws.Range("E21").Value = 8
$ws.Range("F21").Value = "Meetings, working on the poster and panic button test unit"

# Remove the fill from the E/F/D header-ish columns (border-only style) - this happens
# automatically as a side effect of the cellXfs reduction caused by removing unused fill.
# Instead directly clear interior fill (if any) on the D column cells (rows 7-22)
$ws.Range("D7:D22").Interior.Pattern = -4142  # xlPatternNone

# Adjust the sheet view: scroll position, zoom, and selection
$ws.Application.ActiveWindow.Zoom = 87
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F15").Select()
